# Apply the "Add files via upload" edits to the UML Tools comparison sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New data in existing rows -------------------------------------------------

# EclipseUML (row 6) now lists "UML2" under the Files column (G).
$ws.Range("G6").Value = "UML2"

# Sparx Enterprise Architect (row 25): Files cell text changes, and a new
# Associated cell is added.
$ws.Range("G25").Value = "EAP & images"
$ws.Range("H25").Value = "UML2 import"

# Creately (row 26): new Files / Associated cells.
$ws.Range("G26").Value = "CDML"
$ws.Range("H26").Value = "Visio imports"

# Gliffy (row 30): new Files / Associated cells.
$ws.Range("G30").Value = "Proprietary"
$ws.Range("H30").Value = "Visio imports"

# Trace Modeler (row 31): new Cost cell and a wrapped Notes cell.
$ws.Range("E31").Value = "Free"
$ws.Range("I31").Value = "Limited to sequence diagrams"
$ws.Range("I31").WrapText = $true

# --- Brand new rows at the bottom of the table ---------------------------------

$ws.Range("B46").Value = "Stack overflow question on point"
$ws.Range("C46").Value = "Link"
$ws.Hyperlinks.Add($ws.Range("C46"), "https://stackoverflow.com/questions/8384004/eclipse-uml2-tutorial") | Out-Null
$ws.Range("C3").Copy()
$ws.Range("C46").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B47").Value = "Eclipse Foundation"
$ws.Range("C47").Value = "Link"
$ws.Hyperlinks.Add($ws.Range("C47"), "https://www.eclipse.org/modeling/mdt/?project=uml2") | Out-Null
$ws.Range("C3").Copy()
$ws.Range("C47").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("G47").Value = "UML2"

# Leave the selection where the author left it before saving.
$ws.Activate()
$ws.Range("H6").Select()

$wb.Save()
